$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, its column letter, and the new text value.
# Column D holds price strings that can look numeric to Excel's type
# inference (e.g. "1.00", "30.20", "7.90"); force Text via NumberFormat
# "@" while assigning, then restore the default "Normal" style so no
# extra per-cell formatting is left behind (column E values - the
# " +x.xx% " volume strings - are never ambiguous, so they're set
# directly).
$updates = @(
    @{ Cell = 'D2'; Col = 'D'; Value = '66.186.89' },
    @{ Cell = 'E2'; Col = 'E'; Value = '  +1.54%  ' },
    @{ Cell = 'D3'; Col = 'D'; Value = '3.565.28' },
    @{ Cell = 'E3'; Col = 'E'; Value = '  +5.74%  ' },
    @{ Cell = 'D4'; Col = 'D'; Value = '1.00' },
    @{ Cell = 'E4'; Col = 'E'; Value = '  +0.06%  ' },
    @{ Cell = 'D5'; Col = 'D'; Value = '607.48' },
    @{ Cell = 'E5'; Col = 'E'; Value = '  +2.58%  ' },
    @{ Cell = 'D6'; Col = 'D'; Value = '145.14' },
    @{ Cell = 'E6'; Col = 'E'; Value = '  +3.03%  ' },
    @{ Cell = 'D7'; Col = 'D'; Value = '3.564.72' },
    @{ Cell = 'E7'; Col = 'E'; Value = '  +5.75%  ' },
    @{ Cell = 'E8'; Col = 'E'; Value = '  +0.12%  ' },
    @{ Cell = 'E9'; Col = 'E'; Value = '  +4.22%  ' },
    @{ Cell = 'E10'; Col = 'E'; Value = '  +2.59%  ' },
    @{ Cell = 'D11'; Col = 'D'; Value = '7.99' },
    @{ Cell = 'E11'; Col = 'E'; Value = '  +1.45%  ' },
    @{ Cell = 'D12'; Col = 'D'; Value = '0.412' },
    @{ Cell = 'E12'; Col = 'E'; Value = '  +1.96%  ' },
    @{ Cell = 'D13'; Col = 'D'; Value = '4.172.01' },
    @{ Cell = 'E13'; Col = 'E'; Value = '  +5.80%  ' },
    @{ Cell = 'D14'; Col = 'D'; Value = '0.0000208' },
    @{ Cell = 'E14'; Col = 'E'; Value = '  +5.28%  ' },
    @{ Cell = 'D15'; Col = 'D'; Value = '30.20' },
    @{ Cell = 'E15'; Col = 'E'; Value = '  +1.84%  ' },
    @{ Cell = 'D16'; Col = 'D'; Value = '3.565.56' },
    @{ Cell = 'E16'; Col = 'E'; Value = '  +5.69%  ' },
    @{ Cell = 'D17'; Col = 'D'; Value = '66.310.32' },
    @{ Cell = 'E17'; Col = 'E'; Value = '  +1.58%  ' },
    @{ Cell = 'E18'; Col = 'E'; Value = '  -0.65%  ' },
    @{ Cell = 'D19'; Col = 'D'; Value = '11.47' },
    @{ Cell = 'E19'; Col = 'E'; Value = '  +11.43%  ' },
    @{ Cell = 'D20'; Col = 'D'; Value = '6.22' },
    @{ Cell = 'E20'; Col = 'E'; Value = '  +2.44%  ' },
    @{ Cell = 'D21'; Col = 'D'; Value = '14.93' },
    @{ Cell = 'E21'; Col = 'E'; Value = '  +2.15%  ' },
    @{ Cell = 'D22'; Col = 'D'; Value = '431.23' },
    @{ Cell = 'E22'; Col = 'E'; Value = '  +4.30%  ' },
    @{ Cell = 'E23'; Col = 'E'; Value = '  +5.95%  ' },
    @{ Cell = 'D24'; Col = 'D'; Value = '78.69' },
    @{ Cell = 'E24'; Col = 'E'; Value = '  +2.06%  ' },
    @{ Cell = 'D25'; Col = 'D'; Value = '3.707.67' },
    @{ Cell = 'E25'; Col = 'E'; Value = '  +5.72%  ' },
    @{ Cell = 'E26'; Col = 'E'; Value = '  -0.02%  ' },
    @{ Cell = 'D27'; Col = 'D'; Value = '0.0000119' },
    @{ Cell = 'E27'; Col = 'E'; Value = '  +10.26%  ' },
    @{ Cell = 'D28'; Col = 'D'; Value = '2.53' },
    @{ Cell = 'E28'; Col = 'E'; Value = '  +5.21%  ' },
    @{ Cell = 'E29'; Col = 'E'; Value = '  +3.99%  ' },
    @{ Cell = 'D30'; Col = 'D'; Value = '9.17' },
    @{ Cell = 'E30'; Col = 'E'; Value = '  -0.15%  ' },
    @{ Cell = 'D31'; Col = 'D'; Value = '1.00' },
    @{ Cell = 'E31'; Col = 'E'; Value = '  -0.12%  ' },
    @{ Cell = 'E32'; Col = 'E'; Value = '  +2.57%  ' },
    @{ Cell = 'E33'; Col = 'E'; Value = '  -0.61%  ' },
    @{ Cell = 'D34'; Col = 'D'; Value = '3.561.23' },
    @{ Cell = 'E34'; Col = 'E'; Value = '  +5.69%  ' },
    @{ Cell = 'E35'; Col = 'E'; Value = '  +5.20%  ' },
    @{ Cell = 'E36'; Col = 'E'; Value = '  +5.46%  ' },
    @{ Cell = 'E37'; Col = 'E'; Value = '  +0.08%  ' },
    @{ Cell = 'D38'; Col = 'D'; Value = '7.90' },
    @{ Cell = 'E38'; Col = 'E'; Value = '  +5.26%  ' },
    @{ Cell = 'D39'; Col = 'D'; Value = '5.65' },
    @{ Cell = 'E39'; Col = 'E'; Value = '  +3.12%  ' },
    @{ Cell = 'D40'; Col = 'D'; Value = '1.00' },
    @{ Cell = 'E40'; Col = 'E'; Value = '  +0.06%  ' },
    @{ Cell = 'D41'; Col = 'D'; Value = '171.13' },
    @{ Cell = 'E41'; Col = 'E'; Value = '  +2.78%  ' },
    @{ Cell = 'D42'; Col = 'D'; Value = '0.0857' },
    @{ Cell = 'E42'; Col = 'E'; Value = '  +0.57%  ' },
    @{ Cell = 'E43'; Col = 'E'; Value = '  +4.25%  ' },
    @{ Cell = 'D44'; Col = 'D'; Value = '0.896' },
    @{ Cell = 'E44'; Col = 'E'; Value = '  +3.80%  ' },
    @{ Cell = 'D45'; Col = 'D'; Value = '1.94' },
    @{ Cell = 'E45'; Col = 'E'; Value = '  +1.95%  ' },
    @{ Cell = 'D46'; Col = 'D'; Value = '46.04' },
    @{ Cell = 'E46'; Col = 'E'; Value = '  +1.77%  ' },
    @{ Cell = 'D47'; Col = 'D'; Value = '1.21' },
    @{ Cell = 'E47'; Col = 'E'; Value = '  +4.26%  ' },
    @{ Cell = 'D48'; Col = 'D'; Value = '25.94' },
    @{ Cell = 'E48'; Col = 'E'; Value = '  -1.52%  ' },
    @{ Cell = 'D49'; Col = 'D'; Value = '2.38' },
    @{ Cell = 'E49'; Col = 'E'; Value = '  +5.52%  ' },
    @{ Cell = 'E50'; Col = 'E'; Value = '  +1.64%  ' },
    @{ Cell = 'D51'; Col = 'D'; Value = '23.48' },
    @{ Cell = 'E51'; Col = 'E'; Value = '  +16.23%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Col -eq 'D') {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
